$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing header cell (bold, centered, thin border)
# onto the three new header cells so they match the rest of row 1.
$ws.Range("AB1").Copy($ws.Range("AC1:AE1"))

# New header labels for the team-record columns.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Every data row (2-39) gets the same team record: 53 wins, 61 losses, 1 tie.
for ($r = 2; $r -le 39; $r++) {
    $ws.Range("AC$r").Value = 53
    $ws.Range("AD$r").Value = 61
    $ws.Range("AE$r").Value = 1
}
